$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1462.25
$ws.Cells.Item(12, 9).Value = 425
$ws.Cells.Item(12, 11).Value = 425
$ws.Cells.Item(12, 13).Value = -255
$ws.Cells.Item(33, 8).Value = 895.1875
$ws.Cells.Item(33, 9).Value = 1018.8333
$ws.Cells.Item(33, 11).Value = 1018.8333
$ws.Cells.Item(33, 13).Value = -789.8333
$ws.Cells.Item(113, 8).Value = 100016000
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 2995.1428
$ws.Cells.Item(132, 9).Value = 2956.513
$ws.Cells.Item(132, 11).Value = 8869.539000000001
$ws.Cells.Item(132, 13).Value = -6339.539000000001
$ws.Cells.Item(137, 8).Value = 3534.75
$ws.Cells.Item(137, 9).Value = 4760.2
$ws.Cells.Item(137, 11).Value = 14280.6
$ws.Cells.Item(137, 13).Value = -11730.6
$ws.Cells.Item(138, 8).Value = 4832.625
$ws.Cells.Item(138, 9).Value = 949.64
$ws.Cells.Item(138, 10).Value = 9053.261
$ws.Cells.Item(138, 11).Value = 2848.92
$ws.Cells.Item(138, 12).Value = 27159.783
$ws.Cells.Item(138, 13).Value = 2291.08
$ws.Cells.Item(138, 14).Value = -37439.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 5049.467
$ws.Cells.Item(45, 9).Value = 1765.3334
$ws.Cells.Item(45, 10).Value = 9975.666999999999
$ws.Cells.Item(45, 11).Value = 1765.3334
$ws.Cells.Item(45, 12).Value = 9975.666999999999
$ws.Cells.Item(45, 13).Value = -1388.3334
$ws.Cells.Item(45, 14).Value = -10729.667
$ws.Cells.Item(61, 8).Value = 8209.299999999999
$ws.Cells.Item(61, 9).Value = 3619.6
$ws.Cells.Item(61, 11).Value = 3619.6
$ws.Cells.Item(61, 13).Value = -3407.6
$ws.Cells.Item(74, 8).Value = 66428.2
$ws.Cells.Item(74, 9).Value = 95497.82000000001
$ws.Cells.Item(74, 11).Value = 95497.82000000001
$ws.Cells.Item(74, 13).Value = -94623.82000000001
$ws.Cells.Item(77, 8).Value = 66428.2
$ws.Cells.Item(77, 9).Value = 95497.82000000001
$ws.Cells.Item(77, 11).Value = 477489.1
$ws.Cells.Item(77, 13).Value = -473121.1
$ws.Cells.Item(102, 8).Value = 3840.2856
$ws.Cells.Item(102, 9).Value = 3296.6
$ws.Cells.Item(102, 10).Value = 5199.5
$ws.Cells.Item(102, 11).Value = 3296.6
$ws.Cells.Item(102, 12).Value = 5199.5
$ws.Cells.Item(102, 13).Value = -1674.6
$ws.Cells.Item(102, 14).Value = -8443.5
$ws.Cells.Item(122, 8).Value = 21782.084
$ws.Cells.Item(122, 9).Value = 28923.875
$ws.Cells.Item(122, 11).Value = 86771.625
$ws.Cells.Item(122, 13).Value = -84321.625
$ws.Cells.Item(136, 8).Value = 8209.299999999999
$ws.Cells.Item(136, 9).Value = 3619.6
$ws.Cells.Item(136, 11).Value = 10858.8
$ws.Cells.Item(136, 13).Value = -8308.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(96, 8).Value = 26324.715
$ws.Cells.Item(96, 9).Value = 15384.4
$ws.Cells.Item(96, 11).Value = 15384.4
$ws.Cells.Item(96, 13).Value = -12638.4
$ws.Cells.Item(105, 8).Value = 3342.8096
$ws.Cells.Item(105, 9).Value = 2641.1667
$ws.Cells.Item(105, 11).Value = 2641.1667
$ws.Cells.Item(105, 13).Value = -894.1667000000002
$ws.Cells.Item(108, 8).Value = 59380
$ws.Cells.Item(108, 10).Value = 59380
$ws.Cells.Item(108, 12).Value = 59380
$ws.Cells.Item(108, 14).Value = -67060

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 11633380
$ws.Cells.Item(58, 9).Value = 20834934
$ws.Cells.Item(58, 10).Value = 10364.421
$ws.Cells.Item(58, 11).Value = 20834934
$ws.Cells.Item(58, 12).Value = 10364.421
$ws.Cells.Item(58, 13).Value = -20834731
$ws.Cells.Item(58, 14).Value = -10770.421
$ws.Cells.Item(62, 8).Value = 7706
$ws.Cells.Item(62, 9).Value = 7323.6665
$ws.Cells.Item(62, 11).Value = 7323.6665
$ws.Cells.Item(62, 13).Value = -6699.6665
$ws.Cells.Item(65, 8).Value = 7706
$ws.Cells.Item(65, 9).Value = 7323.6665
$ws.Cells.Item(65, 11).Value = 36618.3325
$ws.Cells.Item(65, 13).Value = -33498.3325
$ws.Cells.Item(68, 8).Value = 79995.336
$ws.Cells.Item(68, 10).Value = 79995.336
$ws.Cells.Item(68, 12).Value = 79995.336
$ws.Cells.Item(68, 14).Value = -81493.336
$ws.Cells.Item(71, 8).Value = 79995.336
$ws.Cells.Item(71, 10).Value = 79995.336
$ws.Cells.Item(71, 12).Value = 239986.008
$ws.Cells.Item(71, 14).Value = -247474.008
$ws.Cells.Item(74, 8).Value = 333400000
$ws.Cells.Item(74, 10).Value = 100000
$ws.Cells.Item(74, 12).Value = 100000
$ws.Cells.Item(74, 14).Value = -101748
$ws.Cells.Item(77, 8).Value = 333400000
$ws.Cells.Item(77, 10).Value = 100000
$ws.Cells.Item(77, 12).Value = 300000
$ws.Cells.Item(77, 14).Value = -308736
$ws.Cells.Item(107, 8).Value = 2492.3333
$ws.Cells.Item(107, 9).Value = 2330
$ws.Cells.Item(107, 10).Value = 2573.5
$ws.Cells.Item(107, 11).Value = 2330
$ws.Cells.Item(107, 12).Value = 2573.5
$ws.Cells.Item(107, 13).Value = -410
$ws.Cells.Item(107, 14).Value = -6413.5
$ws.Cells.Item(132, 8).Value = 5320.615
$ws.Cells.Item(132, 9).Value = 2140.0625
$ws.Cells.Item(132, 10).Value = 10409.5
$ws.Cells.Item(132, 11).Value = 6420.1875
$ws.Cells.Item(132, 12).Value = 31228.5
$ws.Cells.Item(132, 13).Value = -3890.1875
$ws.Cells.Item(132, 14).Value = -36288.5
$ws.Cells.Item(134, 8).Value = 7777.5
$ws.Cells.Item(134, 9).Value = 7151.76
$ws.Cells.Item(134, 11).Value = 21455.28
$ws.Cells.Item(134, 13).Value = -18920.28
$ws.Cells.Item(136, 8).Value = 11633380
$ws.Cells.Item(136, 9).Value = 20834934
$ws.Cells.Item(136, 10).Value = 10364.421
$ws.Cells.Item(136, 11).Value = 62504802
$ws.Cells.Item(136, 12).Value = 31093.263
$ws.Cells.Item(136, 13).Value = -62502252
$ws.Cells.Item(136, 14).Value = -36193.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1574.7222
$ws.Cells.Item(131, 10).Value = 4149.25
$ws.Cells.Item(131, 12).Value = 12447.75
$ws.Cells.Item(131, 14).Value = -22527.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 7142.7144
$ws.Cells.Item(21, 9).Value = 4999.75
$ws.Cells.Item(21, 11).Value = 4999.75
$ws.Cells.Item(21, 13).Value = -4826.75
$ws.Cells.Item(30, 8).Value = 7142.7144
$ws.Cells.Item(30, 9).Value = 4999.75
$ws.Cells.Item(30, 11).Value = 4999.75
$ws.Cells.Item(30, 13).Value = -4894.75
$ws.Cells.Item(102, 8).Value = 4281.1333
$ws.Cells.Item(102, 9).Value = 4287.2
$ws.Cells.Item(102, 11).Value = 4287.2
$ws.Cells.Item(102, 13).Value = -2665.2
$ws.Cells.Item(122, 8).Value = 43829.88
$ws.Cells.Item(122, 9).Value = 65698.31
$ws.Cells.Item(122, 10).Value = 4952.6665
$ws.Cells.Item(122, 11).Value = 197094.93
$ws.Cells.Item(122, 12).Value = 14857.9995
$ws.Cells.Item(122, 13).Value = -194644.93
$ws.Cells.Item(122, 14).Value = -19757.9995
$ws.Cells.Item(126, 8).Value = 2551
$ws.Cells.Item(126, 9).Value = 2496.875
$ws.Cells.Item(126, 10).Value = 2599.111
$ws.Cells.Item(126, 11).Value = 7490.625
$ws.Cells.Item(126, 12).Value = 7797.333
$ws.Cells.Item(126, 13).Value = -5020.625
$ws.Cells.Item(126, 14).Value = -12737.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2624.125
$ws.Cells.Item(68, 9).Value = 2248.75
$ws.Cells.Item(68, 10).Value = 2999.5
$ws.Cells.Item(68, 11).Value = 2248.75
$ws.Cells.Item(68, 12).Value = 2999.5
$ws.Cells.Item(68, 13).Value = -1499.75
$ws.Cells.Item(68, 14).Value = -4497.5
$ws.Cells.Item(71, 8).Value = 2624.125
$ws.Cells.Item(71, 9).Value = 2248.75
$ws.Cells.Item(71, 10).Value = 2999.5
$ws.Cells.Item(71, 11).Value = 11243.75
$ws.Cells.Item(71, 12).Value = 14997.5
$ws.Cells.Item(71, 13).Value = -7499.75
$ws.Cells.Item(71, 14).Value = -22485.5
$ws.Cells.Item(92, 8).Value = 46636.5
$ws.Cells.Item(92, 10).Value = 46636.5
$ws.Cells.Item(92, 12).Value = 46636.5
$ws.Cells.Item(92, 14).Value = -51628.5
$ws.Cells.Item(100, 8).Value = 3351.375
$ws.Cells.Item(100, 10).Value = 3628
$ws.Cells.Item(100, 12).Value = 3628
$ws.Cells.Item(100, 14).Value = -4710
$ws.Cells.Item(136, 8).Value = 7099
$ws.Cells.Item(136, 9).Value = 3414.44
$ws.Cells.Item(136, 10).Value = 11947.105
$ws.Cells.Item(136, 11).Value = 10243.32
$ws.Cells.Item(136, 12).Value = 35841.315
$ws.Cells.Item(136, 13).Value = -7693.32
$ws.Cells.Item(136, 14).Value = -40941.315

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 54794
$ws.Cells.Item(103, 10).Value = 54794
$ws.Cells.Item(103, 12).Value = 54794
$ws.Cells.Item(103, 14).Value = -57138
$ws.Cells.Item(112, 8).Value = 39994.5
$ws.Cells.Item(112, 10).Value = 39994.5
$ws.Cells.Item(112, 12).Value = 39994.5
$ws.Cells.Item(112, 14).Value = -42948.5
$ws.Cells.Item(122, 8).Value = 6187.3125
$ws.Cells.Item(122, 9).Value = 4399
$ws.Cells.Item(122, 11).Value = 13197
$ws.Cells.Item(122, 13).Value = -10747
$ws.Cells.Item(132, 8).Value = 18551786
$ws.Cells.Item(132, 9).Value = 20839508
$ws.Cells.Item(132, 11).Value = 62518524
$ws.Cells.Item(132, 13).Value = -62515994
$ws.Cells.Item(136, 8).Value = 66672784
$ws.Cells.Item(136, 9).Value = 250001740
$ws.Cells.Item(136, 11).Value = 750005220
$ws.Cells.Item(136, 13).Value = -750002670
